$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "66.599.35"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "3.593.87"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.21"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.29"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.04"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "4.204.90"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.04"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.590.37"
$ws.Range("D16").Value = "66.671.74"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.49"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.07"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.73"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.621"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.17"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "3.736.79"
$ws.Range("E24").Value = "  +0.89%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.25"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "3.592.52"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.55"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.85"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "177.01"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0858"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.898"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +8.02%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.04"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.12"
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -1.42%  "
